# Generate Report for Handoff
# Updates the localization-status workbook: status moves from
# "In Translation" to "Ready for handoff" and the associated timestamps
# are refreshed; the Status/date columns are also widened.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---------------------------------------------------
# E2 = zh-cn status, F2 = de-de status, G2 = Latest HO Xliff Generate Date
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-23 02:57:11"

# --- zh-cn sheet --------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-23 02:57:03"

# --- de-de sheet --------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-23 02:57:11"

# --- Widen the Status/date columns --------------------------------------
# Target OOXML column width is 17.2159881591797 characters; the closest
# value reachable through the ColumnWidth property (quantized in 1/6
# character steps by this engine) is 17.1666... , reached by requesting
# 16.3333... characters.
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333332
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333332
